$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '94.910.55'
$ws.Range("E2").Value = '  -2.20%  '

$ws.Range("D3").Value = '3.559.34'
$ws.Range("E3").Value = '  -1.52%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.15'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.76%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '652.73'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.56%  '

$ws.Range("E7").Value = '  -2.75%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.397'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.47%  '

$ws.Range("E9").Value = '  +0.12%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.997'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.38%  '

$ws.Range("D11").Value = '3.556.70'
$ws.Range("E11").Value = '  -1.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.202'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.29%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.15'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.60%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.41'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.07%  '

$ws.Range("D15").Value = '4.222.31'
$ws.Range("E15").Value = '  -1.86%  '

$ws.Range("D16").Value = '94.881.41'
$ws.Range("E16").Value = '  -2.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000252'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.52%  '

$ws.Range("D18").Value = '3.557.32'
$ws.Range("E18").Value = '  -1.49%  '

$ws.Range("E19").Value = '  -4.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.65'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.68'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.60%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.44'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.71%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '506.03'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.477'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -5.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.75'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.57%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000194'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.35%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '94.76'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.77%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.53'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.34%  '

$ws.Range("D29").Value = '3.752.13'
$ws.Range("E29").Value = '  -1.45%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.01'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -4.92%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.142'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -3.09%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.43'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.29%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.10%  '

$ws.Range("E34").Value = '  -0.95%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.176'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.46%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.72'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +3.92%  '

$ws.Range("E37").Value = '  +11.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.553'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.70%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.42'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +6.08%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '581.79'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.49%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.150'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.71%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.901'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.64%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.80'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.60%  '

$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.70'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.62%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '34.47'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +31.73%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.26'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.85%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.38'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.86%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0411'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -6.27%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.55'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.26%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.07'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.15%  '

